# Update the carjacking arrests workbook with data through 2021-12-15

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (sheet tab name AND workbook.xml <sheet name="...">)
$ws.Name = "Through 2021-12-15"

# Row 13 (November) - only T,U,V changed
$ws.Range("T13").Value = 6
$ws.Range("U13").Value = 194
$ws.Range("V13").Value = 0.03

# Row 14 (December) - label + many cells changed
$ws.Range("A14").Value = "December (through 12-15)"
$ws.Range("C14").Value = 13
$ws.Range("D14").Value = 0.1875
$ws.Range("F14").Value = 41
$ws.Range("G14").Value = 0.1087
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = 47
$ws.Range("J14").Value = 0.0962
$ws.Range("L14").Value = 31
$ws.Range("M14").Value = 0.0882
$ws.Range("O14").Value = 22
$ws.Range("P14").Value = 0.12
$ws.Range("Q14").Value = 4
$ws.Range("R14").Value = 73
$ws.Range("S14").Value = 0.0519
$ws.Range("U14").Value = 111
$ws.Range("V14").Value = 0.0089

# Row 15 (Total) - sums updated
$ws.Range("C15").Value = 271
$ws.Range("D15").Value = 0.1173
$ws.Range("F15").Value = 544
$ws.Range("G15").Value = 0.1067
$ws.Range("H15").Value = 68
$ws.Range("I15").Value = 805
$ws.Range("J15").Value = 0.0779
$ws.Range("L15").Value = 639
$ws.Range("M15").Value = 0.1075
$ws.Range("O15").Value = 502
$ws.Range("P15").Value = 0.102
$ws.Range("Q15").Value = 68
$ws.Range("R15").Value = 1273
$ws.Range("S15").Value = 0.0507
$ws.Range("T15").Value = 102
$ws.Range("U15").Value = 1653
$ws.Range("V15").Value = 0.0581
